$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 122, shifting existing rows 122-204 down to 123-205
$ws.Rows("122:122").Insert()

# Populate the new row 122 with the new weekly data point
$ws.Range("A122").Value = 10
$ws.Range("B122").Value = "Vega Modelo de Temuco"
$ws.Range("C122").Value = "La Araucanía"
$ws.Range("D122").Value = 44603
$ws.Range("E122").Value = 9
$ws.Range("F122").Value = 100112052
$ws.Range("G122").Value = "Albahaca"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 35
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 6000
$ws.Range("N122").Value = "$/paquete"
$ws.Range("O122").Value = "Región de La Araucanía"
$ws.Range("P122").Value = 6000
$ws.Range("Q122").Value = 1
$ws.Range("R122").Value = "Hortaliza"
